# The roster previously led with two "term" columns (term_code, term_name)
# that are no longer wanted; drop them so the sheet starts at person_name.
# Deleting the columns shifts everything else two columns to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:B1").EntireColumn.Delete()

# Add a new roster entry for Kian, an Alumni (name + role only).
$ws.Range("A10").Value = "Kian"
$ws.Range("E10").Value = "Alumni"
